$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column E / F "Razon social" / "Nombre Fantasia" punctuation fixes ---
# Comma separators between multiple names replaced with periods (per source diff).
$ws.Cells.Item(93, 5).Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Cells.Item(131, 5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Cells.Item(210, 5).Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Cells.Item(209, 5).Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Cells.Item(222, 5).Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Cells.Item(277, 6).Value = "CLERICE. MIGUEL ANGEL"

# --- Column H "Importe" reformat: "1.234,56" (es-AR) -> "1234.56" (plain) ---
# Values must stay TEXT (shared-string) cells, matching the original cell type,
# so we build them as formulas yielding a text string, then convert the whole
# range to static values via Copy / PasteSpecial (xlPasteValues) - this avoids
# Excel re-interpreting "28000.00" as a number (which would also add a new
# number-format style to the cell).
$newImportes = @(
  "28000.00",
  "211200.00",
  "36000.00",
  "184800.00",
  "306000.00",
  "14600.00",
  "1375500.00",
  "1681000.00",
  "1443.00",
  "1868.00",
  "1428.85",
  "2525.00",
  "4840.00",
  "1529144.00",
  "12045.00",
  "521200.00",
  "8027.10",
  "3870.00",
  "3600.00",
  "17280.00",
  "154054.00",
  "737780.48",
  "116528.00",
  "418282.66",
  "161970.78",
  "79799.36",
  "13115.00",
  "26249.04",
  "82020.00",
  "2200.00",
  "9315.60",
  "79459.56",
  "21377.00",
  "33499.88",
  "6800.00",
  "21330.00",
  "27200.00",
  "80000.00",
  "5000.00",
  "264.00",
  "406755.48",
  "367876.85",
  "29832.00",
  "700.00",
  "16378.98",
  "4555.16",
  "954699.00",
  "5980.00",
  "77681.95",
  "1576.61",
  "9300.00",
  "18701.39",
  "101064.40",
  "11360.00",
  "4067.00",
  "8752.00",
  "13232.50",
  "22586.89",
  "95156.21",
  "2862.67",
  "13001.00",
  "8925.00",
  "79893.00",
  "11263.90",
  "1700.00",
  "4365.64",
  "668109.00",
  "26278.24",
  "4050.00",
  "48850.00",
  "6000.00",
  "26800.00",
  "232000.00",
  "168880.00",
  "1600.00",
  "18132.00",
  "15640.00",
  "12500.00",
  "7800.00",
  "9800.00",
  "212475.00",
  "410.00",
  "6700.00",
  "26500.00",
  "1900.00",
  "16000.00",
  "25000.00",
  "60000.00",
  "1470.00",
  "28990.00",
  "900.00",
  "7779.00",
  "47360.00",
  "6500.00",
  "3700.00",
  "11135.00",
  "5959.00",
  "34619.45",
  "37.48",
  "58497.91",
  "36558.14",
  "22218.62",
  "11120.00",
  "90.00",
  "493.38",
  "93367.00",
  "15600.00",
  "64888.00",
  "37974.00",
  "15456.22",
  "189338.90",
  "13160.00",
  "1950.00",
  "40657.00",
  "15570.00",
  "2275.00",
  "2100.00",
  "12180.00",
  "1375.00",
  "25340.00",
  "64225.00",
  "26720.00",
  "5077.00",
  "19540.22",
  "9520.00",
  "549.80",
  "2905.00",
  "300.00",
  "1600.00",
  "365.00",
  "264.00",
  "704.00",
  "2760.00",
  "22800.00",
  "270.00",
  "21600.00",
  "1905828.00",
  "2000.00",
  "5000.00",
  "29000.00",
  "30000.00",
  "8000.00",
  "4000.00",
  "10500.00",
  "10000.00",
  "10000.00",
  "16000.00",
  "35000.00",
  "44000.00",
  "59800.00",
  "40000.00",
  "55000.00",
  "36000.00",
  "350.00",
  "67000.00",
  "9084.33",
  "1971.00",
  "21264.00",
  "433.81",
  "1760.00",
  "3876.87",
  "11100.00",
  "70000.00",
  "18000.00",
  "33000.00",
  "24000.00",
  "22000.00",
  "21000.00",
  "12000.00",
  "20000.00",
  "17000.00",
  "6720.00",
  "24000.00",
  "10000.00",
  "20000.00",
  "9000.00",
  "18000.00",
  "3500.00",
  "10000.00",
  "20000.00",
  "10000.00",
  "42000.00",
  "20000.00",
  "20000.00",
  "28500.00",
  "10000.00",
  "41733.20",
  "46000.00",
  "20000.00",
  "24000.00",
  "9000.00",
  "24000.00",
  "28500.00",
  "51000.00",
  "20000.00",
  "10000.00",
  "82557.10",
  "3500.00",
  "74860.00",
  "20000.00",
  "71999.98",
  "12000.00",
  "11500.00",
  "400.00",
  "14000.00",
  "5562.73",
  "14050.00",
  "6655.00",
  "19070.00",
  "21157.00",
  "6482.00",
  "27050.00",
  "2914.00",
  "1468.74",
  "3400.00",
  "95414.20",
  "982.00",
  "9869.07",
  "15312.00",
  "37193.00",
  "24907.98",
  "24025.00",
  "660.00",
  "1463.04",
  "3159.00",
  "3470.00",
  "1265146.41",
  "1080.00",
  "7114.80",
  "80000.00",
  "40000.00",
  "40000.00",
  "40000.00",
  "80000.00",
  "40000.00",
  "55000.00",
  "40000.00",
  "40000.00",
  "80000.00",
  "80000.00",
  "74000.00",
  "6300.00",
  "5227021.55",
  "216000.00",
  "240000.00",
  "2425627.36",
  "30000.00",
  "7472405.19",
  "241000.00",
  "241000.00",
  "250200.00",
  "242200.00",
  "243700.00",
  "245500.00",
  "451000.00",
  "241000.00",
  "573450.00",
  "513000.00",
  "284650.00",
  "241000.00",
  "241000.00",
  "482000.00",
  "370400.00",
  "507800.00",
  "685300.00",
  "451000.00",
  "708200.00",
  "482000.00",
  "247850.00",
  "311587.11",
  "176300.00",
  "262350.00",
  "14000.00",
  "10963063.74",
  "22000.00",
  "49760.00",
  "118000.00",
  "74200.00",
  "92820.00",
  "25540.00",
  "76800.00",
  "10500.00",
  "5500.00",
  "454000.00",
  "11610.00",
  "18500.00",
  "3900.00",
  "27000.00"
)

for ($i = 0; $i -lt $newImportes.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 8).Formula = '="' + $newImportes[$i] + '"'
}

$importeRange = $ws.Range("H2:H289")
$importeRange.Copy()
$importeRange.PasteSpecial(-4163)

